$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = -8.608899999999997
$ws.Range("D6").Value = -8.274599999999996
$ws.Range("C7").Value = -12.5447
$ws.Range("A8").Value = -21.2192
$ws.Range("D9").Value = -8.180199999999997
$ws.Range("A10").Value = -20.43709999999998
$ws.Range("D10").Value = -6.171999999999998
$ws.Range("A12").Value = -22.55410000000003
$ws.Range("B13").Value = 5.7296
$ws.Range("A18").Value = -22.33240000000001
$ws.Range("C20").Value = -14.28849999999999
